# Apply updated Betfair back/lay odds values to Sheet1.
# Values below reflect the cell-level deltas captured in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (FK Radnicki 1923 vs Cukaricki)
$ws.Range("H2").Value = 2.5
$ws.Range("J2").Value = 3.25
$ws.Range("L2").Value = 1.27
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 1.01
$ws.Range("P2").Value = 2.08
$ws.Range("Q2").Value = 1.56
$ws.Range("R2").Value = 1.42
$ws.Range("W2").Value = 1.55

# Row 3 (Young Boys vs Grasshoppers Zurich)
$ws.Range("H3").Value = 5.2
$ws.Range("N3").Value = 6
$ws.Range("T3").Value = 1.52
$ws.Range("W3").Value = 2.56
$ws.Range("Y3").Value = 34
$ws.Range("AN3").Value = 5.8

# Row 4 (FC Zurich vs Lugano)
$ws.Range("F4").Value = 2.74
$ws.Range("L4").Value = 1.28
$ws.Range("N4").Value = 4.5
$ws.Range("O4").Value = 1.23
$ws.Range("P4").Value = 2.2
$ws.Range("Q4").Value = 1.59
$ws.Range("R4").Value = 1.49
$ws.Range("S4").Value = 2.72
$ws.Range("T4").Value = 1.59
$ws.Range("W4").Value = 1.5
$ws.Range("X4").Value = 24
$ws.Range("Y4").Value = 16.5
$ws.Range("AB4").Value = 17.5
$ws.Range("AC4").Value = 11
$ws.Range("AD4").Value = 15
$ws.Range("AI4").Value = 40
$ws.Range("AK4").Value = 36
$ws.Range("AL4").Value = 44
$ws.Range("AM4").Value = 80
$ws.Range("AN4").Value = 24
$ws.Range("AO4").Value = 21

# Row 5 (Luzern vs FC Basel)
$ws.Range("I5").Value = 2.16
$ws.Range("K5").Value = 4.1
$ws.Range("R5").Value = 1.65
$ws.Range("S5").Value = 2.32
$ws.Range("V5").Value = 1.87

# Row 6 (Dundee Utd vs Celtic)
$ws.Range("I6").Value = 1.51
$ws.Range("K6").Value = 5.1
$ws.Range("P6").Value = 2.42
$ws.Range("V6").Value = 2.96
$ws.Range("AC6").Value = 11.5
$ws.Range("AJ6").Value = 220
